$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 44.545456
$ws.Range("I11").Value = 44.545456
$ws.Range("K11").Value = 44.545456
$ws.Range("M11").Value = 95.454544
$ws.Range("H53").Value = 3698.9092
$ws.Range("I53").Value = 250
$ws.Range("J53").Value = 4465.3335
$ws.Range("K53").Value = 250
$ws.Range("L53").Value = 4465.3335
$ws.Range("M53").Value = 387
$ws.Range("N53").Value = -5739.3335
$ws.Range("H113").Value = 83337840
$ws.Range("I113").Value = 500000740
$ws.Range("J113").Value = 5260.3
$ws.Range("K113").Value = 500000740
$ws.Range("L113").Value = 5260.3
$ws.Range("M113").Value = -499997486
$ws.Range("N113").Value = -11768.3
$ws.Range("H129").Value = 345830.56
$ws.Range("I129").Value = 495
$ws.Range("J129").Value = 371410.97
$ws.Range("K129").Value = 1485
$ws.Range("L129").Value = 1114232.91
$ws.Range("M129").Value = 3515
$ws.Range("N129").Value = -1124232.91
$ws.Range("H132").Value = 18520402
$ws.Range("I132").Value = 20001920
$ws.Range("J132").Value = 1426.5
$ws.Range("K132").Value = 60005760
$ws.Range("L132").Value = 4279.5
$ws.Range("M132").Value = -60003230
$ws.Range("N132").Value = -9339.5
$ws.Range("H137").Value = 93345.67999999999
$ws.Range("I137").Value = 110044.87
$ws.Range("K137").Value = 330134.61
$ws.Range("M137").Value = -327584.61
$ws.Range("H138").Value = 3036.4707
$ws.Range("I138").Value = 2610.5293
$ws.Range("J138").Value = 3142.9558
$ws.Range("K138").Value = 7831.5879
$ws.Range("L138").Value = 9428.867400000001
$ws.Range("M138").Value = -2691.5879
$ws.Range("N138").Value = -19708.8674

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11692.091
$ws.Range("I32").Value = 8544.16
$ws.Range("J32").Value = 21529.375
$ws.Range("K32").Value = 8544.16
$ws.Range("L32").Value = 21529.375
$ws.Range("M32").Value = -8257.16
$ws.Range("N32").Value = -22103.375
$ws.Range("H97").Value = 435.42856
$ws.Range("I97").Value = 435.5
$ws.Range("K97").Value = 435.5
$ws.Range("M97").Value = 60.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3646.75
$ws.Range("I20").Value = 3269.8
$ws.Range("J20").Value = 4275
$ws.Range("K20").Value = 3269.8
$ws.Range("L20").Value = 4275
$ws.Range("M20").Value = -3022.8
$ws.Range("N20").Value = -4769
$ws.Range("H52").Value = 28729.834
$ws.Range("J52").Value = 28729.834
$ws.Range("L52").Value = 28729.834
$ws.Range("N52").Value = -29255.834
$ws.Range("H70").Value = 82033.336
$ws.Range("J70").Value = 82033.336
$ws.Range("L70").Value = 82033.336
$ws.Range("N70").Value = -82619.336
$ws.Range("H73").Value = 82033.336
$ws.Range("J73").Value = 82033.336
$ws.Range("L73").Value = 82033.336
$ws.Range("N73").Value = -84061.336
$ws.Range("H99").Value = 1244.2858
$ws.Range("I99").Value = 1320
$ws.Range("J99").Value = 966.6667
$ws.Range("K99").Value = 1320
$ws.Range("L99").Value = 966.6667
$ws.Range("M99").Value = 178
$ws.Range("N99").Value = -3962.6667
$ws.Range("H105").Value = 5418794
$ws.Range("I105").Value = 11112663
$ws.Range("J105").Value = 2002472.5
$ws.Range("K105").Value = 11112663
$ws.Range("L105").Value = 2002472.5
$ws.Range("M105").Value = -11110916
$ws.Range("N105").Value = -2005966.5
$ws.Range("H112").Value = 38999.5
$ws.Range("J112").Value = 38999.5
$ws.Range("L112").Value = 38999.5
$ws.Range("N112").Value = -41953.5
$ws.Range("H121").Value = 28729.834
$ws.Range("J121").Value = 28729.834
$ws.Range("L121").Value = 28729.834
$ws.Range("N121").Value = -32223.834

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 33136.5
$ws.Range("J141").Value = 33136.5
$ws.Range("L141").Value = 33136.5
$ws.Range("N141").Value = -43496.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6583.3335
$ws.Range("I56").Value = 6583.3335
$ws.Range("K56").Value = 6583.3335
$ws.Range("M56").Value = -6053.3335
$ws.Range("H131").Value = 728.05
$ws.Range("J131").Value = 759.7252999999999
$ws.Range("L131").Value = 2279.1759
$ws.Range("N131").Value = -12359.1759

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4050.625
$ws.Range("I80").Value = 3763
$ws.Range("K80").Value = 3763
$ws.Range("M80").Value = -2765
$ws.Range("H83").Value = 4050.625
$ws.Range("I83").Value = 3763
$ws.Range("K83").Value = 18815
$ws.Range("M83").Value = -13823
$ws.Range("H97").Value = 1692.8572
$ws.Range("I97").Value = 667.8
$ws.Range("J97").Value = 4255.5
$ws.Range("K97").Value = 667.8
$ws.Range("L97").Value = 4255.5
$ws.Range("M97").Value = -171.8
$ws.Range("N97").Value = -5247.5
$ws.Range("H132").Value = 3753831.8
$ws.Range("I132").Value = 5776495
$ws.Range("K132").Value = 17329485
$ws.Range("M132").Value = -17326955

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4066.6667
$ws.Range("I22").Value = 4066.6667
$ws.Range("K22").Value = 4066.6667
$ws.Range("M22").Value = -3771.6667
$ws.Range("H27").Value = 4066.6667
$ws.Range("I27").Value = 4066.6667
$ws.Range("K27").Value = 4066.6667
$ws.Range("M27").Value = -3959.6667
$ws.Range("H40").Value = 62120.4
$ws.Range("I40").Value = 92633.08
$ws.Range("J40").Value = 5454
$ws.Range("K40").Value = 92633.08
$ws.Range("L40").Value = 5454
$ws.Range("M40").Value = -92497.08
$ws.Range("N40").Value = -5726
$ws.Range("H82").Value = 2564.0908
$ws.Range("J82").Value = 3250.75
$ws.Range("L82").Value = 3250.75
$ws.Range("N82").Value = -3972.75
$ws.Range("H85").Value = 2564.0908
$ws.Range("J85").Value = 3250.75
$ws.Range("L85").Value = 3250.75
$ws.Range("N85").Value = -5746.75
$ws.Range("H100").Value = 2118.4119
$ws.Range("I100").Value = 2101.25
$ws.Range("J100").Value = 2123.6924
$ws.Range("K100").Value = 2101.25
$ws.Range("L100").Value = 2123.6924
$ws.Range("M100").Value = -1560.25
$ws.Range("N100").Value = -3205.6924
$ws.Range("H122").Value = 1965324.4
$ws.Range("I122").Value = 2804892
$ws.Range("K122").Value = 8414676
$ws.Range("M122").Value = -8412226
$ws.Range("H132").Value = 3247.4138
$ws.Range("I132").Value = 2576.1365
$ws.Range("K132").Value = 7728.4095
$ws.Range("M132").Value = -5198.4095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 32261726
$ws.Range("I136").Value = 50001980
$ws.Range("J136").Value = 6718.5454
$ws.Range("K136").Value = 150005940
$ws.Range("L136").Value = 20155.6362
$ws.Range("M136").Value = -150003390
$ws.Range("N136").Value = -25255.6362

Write-Host "Applied all changes"